# Weekly fruit/vegetable data refresh for "Vega Monumental Concepción - Limón".
# Two new price observations were recorded, so two new rows are inserted
# right before the existing row 504 (pushing the former rows 504-550 down
# to 506-552, same as Excel's normal "insert rows" behaviour).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 504-505; everything from the old row 504 onward
# shifts down by two rows (old 504 -> 506, ..., old 550 -> 552).
$ws.Range("A504:A505").EntireRow.Insert()

# Populate the two newly inserted rows with the new observations
# (columns A..T, same layout as every other row in the sheet).
$row504 = @(11, "Vega Monumental Concepción", "Bíobío", 44783, 8, "Fruta", 100102, "Cítricos", 100102003, "Limón", "Sin especificar", "1a amarillo", 450, 5000, 5500, 5278, "`$/malla 16 kilos", "Provincia de Melipilla", 330, 16)
$row505 = @(11, "Vega Monumental Concepción", "Bíobío", 44783, 8, "Fruta", 100102, "Cítricos", 100102003, "Limón", "Sin especificar", "2a amarillo", 350, 4000, 4500, 4286, "`$/malla 16 kilos", "Provincia de Melipilla", 268, 16)

for ($col = 1; $col -le $row504.Length; $col++) {
    $ws.Cells.Item(504, $col).Value = $row504[$col - 1]
}

for ($col = 1; $col -le $row505.Length; $col++) {
    $ws.Cells.Item(505, $col).Value = $row505[$col - 1]
}
